$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear old data rows (A2:H21) so orphaned shared strings are dropped on save
$ws.Range("A2:H21").ClearContents()

# Shrink from 20 data rows to 18 data rows (rows 20 and 21 go away entirely)
$ws.Range("A20:H21").EntireRow.Delete()

# Column B (Dia / dates) - write top-to-bottom first so shared-string order matches
$c = $ws.Range("B2"); $c.Value = "'2025-05-14"; $c.Style = "Normal"
$c = $ws.Range("B3"); $c.Value = "'2025-05-14"; $c.Style = "Normal"
$c = $ws.Range("B4"); $c.Value = "'2025-05-14"; $c.Style = "Normal"
$c = $ws.Range("B5"); $c.Value = "'2025-05-14"; $c.Style = "Normal"
$c = $ws.Range("B6"); $c.Value = "'2025-05-14"; $c.Style = "Normal"
$c = $ws.Range("B7"); $c.Value = "'2025-05-15"; $c.Style = "Normal"
$c = $ws.Range("B8"); $c.Value = "'2025-05-16"; $c.Style = "Normal"
$c = $ws.Range("B9"); $c.Value = "'2025-05-19"; $c.Style = "Normal"
$c = $ws.Range("B10"); $c.Value = "'2025-05-19"; $c.Style = "Normal"
$c = $ws.Range("B11"); $c.Value = "'2025-05-19"; $c.Style = "Normal"
$c = $ws.Range("B12"); $c.Value = "'2025-05-19"; $c.Style = "Normal"
$c = $ws.Range("B13"); $c.Value = "'2025-05-19"; $c.Style = "Normal"
$c = $ws.Range("B14"); $c.Value = "'2025-05-20"; $c.Style = "Normal"
$c = $ws.Range("B15"); $c.Value = "'2025-05-21"; $c.Style = "Normal"
$c = $ws.Range("B16"); $c.Value = "'2025-05-21"; $c.Style = "Normal"
$c = $ws.Range("B17"); $c.Value = "'2025-05-22"; $c.Style = "Normal"
$c = $ws.Range("B18"); $c.Value = "'2025-05-22"; $c.Style = "Normal"
$c = $ws.Range("B19"); $c.Value = "'2025-05-22"; $c.Style = "Normal"

# Column D (cliente)
$ws.Range("D2").Value = "CONDOMINIO CRISTAL TOWER"
$ws.Range("D3").Value = "MUSASHI DA AMAZONIA LTDA"
$ws.Range("D4").Value = "MUSASHI DA AMAZONIA LTDA"
$ws.Range("D5").Value = "V V REFEICOES LTDA"
$ws.Range("D6").Value = "MUSASHI DA AMAZONIA LTDA"
$ws.Range("D7").Value = "SAT BRAS INDUSTRIA ELETRONICA DA AMAZONIA LTDA."
$ws.Range("D8").Value = "SERVICO NACIONAL DE APRENDIZAGEM COMERCIAL"
$ws.Range("D9").Value = "NORTH HOTEIS EIRELI"
$ws.Range("D10").Value = "TEL TELECOMUNICACOES LTDA."
$ws.Range("D11").Value = "TEL TELECOMUNICACOES LTDA."
$ws.Range("D12").Value = "BEIRA ALTA INDUSTRIAL LTDA"
$ws.Range("D13").Value = "DSB COMERCIO DE ALIMENTOS EIRELI"
$ws.Range("D14").Value = "L T MARTINS LTDA"
$ws.Range("D15").Value = "METALURGICA SATO DA AMAZONIA LTDA"
$ws.Range("D16").Value = "V V REFEICOES LTDA"
$ws.Range("D17").Value = "RODRIGUES INDUSTRIA E COMERCIO DE COLCHOES LTDA"
$ws.Range("D18").Value = "AMAZONIA REFEICOES E SERVICOS LTDA"
$ws.Range("D19").Value = "AMAZONIA REFEICOES E SERVICOS LTDA"

# Column E (id_produto) - force text so leading zeros survive
$c = $ws.Range("E2"); $c.Value = "'000084"; $c.Style = "Normal"
$c = $ws.Range("E3"); $c.Value = "'000307"; $c.Style = "Normal"
$c = $ws.Range("E4"); $c.Value = "'000415"; $c.Style = "Normal"
$c = $ws.Range("E5"); $c.Value = "'000208"; $c.Style = "Normal"
$c = $ws.Range("E6"); $c.Value = "'000842"; $c.Style = "Normal"
$c = $ws.Range("E7"); $c.Value = "'010072"; $c.Style = "Normal"
$c = $ws.Range("E8"); $c.Value = "'000328"; $c.Style = "Normal"
$c = $ws.Range("E9"); $c.Value = "'010040"; $c.Style = "Normal"
$c = $ws.Range("E10"); $c.Value = "'000032"; $c.Style = "Normal"
$c = $ws.Range("E11"); $c.Value = "'000034"; $c.Style = "Normal"
$c = $ws.Range("E12"); $c.Value = "'000083"; $c.Style = "Normal"
$c = $ws.Range("E13"); $c.Value = "'000717"; $c.Style = "Normal"
$c = $ws.Range("E14"); $c.Value = "'000098"; $c.Style = "Normal"
$c = $ws.Range("E15"); $c.Value = "'000032"; $c.Style = "Normal"
$c = $ws.Range("E16"); $c.Value = "'000877"; $c.Style = "Normal"
$c = $ws.Range("E17"); $c.Value = "'000897"; $c.Style = "Normal"
$c = $ws.Range("E18"); $c.Value = "'000999"; $c.Style = "Normal"
$c = $ws.Range("E19"); $c.Value = "'000426"; $c.Style = "Normal"

# Column F (produto)
$ws.Range("F2").Value = "SACO DE LIXO 200L PRETO 0.10 REFORCADO - PCT C/100 UND"
$ws.Range("F3").Value = "DISCO REMOVEDOR PRETO P ENCERADEIRA 510MM"
$ws.Range("F4").Value = "DETERGENTE DESENGRAX MAX PINE AUDAX 5L"
$ws.Range("F5").Value = "MASCARA DESCARTAVEL TRIPLA TALGE C/50 UND"
$ws.Range("F6").Value = "SACO DE LIXO 200L COMUM PACOTINHO C/5 UND"
$ws.Range("F7").Value = "FITA DUREX 48 X 50 TRANSPARENTE"
$ws.Range("F8").Value = "REFIL MOP UMIDO 290G PONTA DOBRADA C LOOP E CINTA"
$ws.Range("F9").Value = "LAPIS PRETO N2 UND LEOELEO"
$ws.Range("F10").Value = "LIMPADOR VEJA MULTIUSO GOLD 500ML"
$ws.Range("F11").Value = "LUVA DE LIMPEZA LONGA NITRILICA VERDE TALGE TAM M PAR"
$ws.Range("F12").Value = "SACO DE LIXO 50L PRETO COMUM - PCT C/100 UND"
$ws.Range("F13").Value = "DETERGENTE NEUTRO BRINORT 500ML"
$ws.Range("F14").Value = "PANO DE CHAO FLANELADO C REFORCADO ITATEX 42x62CM"
$ws.Range("F15").Value = "LIMPADOR VEJA MULTIUSO GOLD 500ML"
$ws.Range("F16").Value = "CABO DE ALUMINIO NOBRE 140 CM COM PONTEIRA"
$ws.Range("F17").Value = "AGUA SANITARIA GLOBO SAN 5L"
$ws.Range("F18").Value = "ESCOVA DE ACO 3 FILEIRAS"
$ws.Range("F19").Value = "KIT DESCASCADOR DE LEGUMES KEITA"

# Numeric / boolean columns A, C, G, H
$ws.Range("A2").Value = 0; $ws.Range("C2").Value = 30; $ws.Range("G2").Value = 19; $ws.Range("H2").Value = $false
$ws.Range("A3").Value = 3; $ws.Range("C3").Value = 160; $ws.Range("G3").Value = 62; $ws.Range("H3").Value = $false
$ws.Range("A4").Value = 11; $ws.Range("C4").Value = 200; $ws.Range("G4").Value = 161; $ws.Range("H4").Value = $true
$ws.Range("A5").Value = 16; $ws.Range("C5").Value = 100; $ws.Range("G5").Value = 311; $ws.Range("H5").Value = $false
$ws.Range("A6").Value = 7; $ws.Range("C6").Value = 400; $ws.Range("G6").Value = 306; $ws.Range("H6").Value = $false
$ws.Range("A7").Value = 5; $ws.Range("C7").Value = 80; $ws.Range("G7").Value = 803; $ws.Range("H7").Value = $false
$ws.Range("A8").Value = 14; $ws.Range("C8").Value = 100; $ws.Range("G8").Value = 85; $ws.Range("H8").Value = $false
$ws.Range("A9").Value = 15; $ws.Range("C9").Value = 400; $ws.Range("G9").Value = 204; $ws.Range("H9").Value = $false
$ws.Range("A10").Value = 8; $ws.Range("C10").Value = 120; $ws.Range("G10").Value = 1200; $ws.Range("H10").Value = $false
$ws.Range("A11").Value = 4; $ws.Range("C11").Value = 60; $ws.Range("G11").Value = 85; $ws.Range("H11").Value = $false
$ws.Range("A12").Value = 2; $ws.Range("C12").Value = 60; $ws.Range("G12").Value = 56; $ws.Range("H12").Value = $false
$ws.Range("A13").Value = 17; $ws.Range("C13").Value = 504; $ws.Range("G13").Value = 1784; $ws.Range("H13").Value = $false
$ws.Range("A14").Value = 1; $ws.Range("C14").Value = 400; $ws.Range("G14").Value = 2186; $ws.Range("H14").Value = $false
$ws.Range("A15").Value = 9; $ws.Range("C15").Value = 84; $ws.Range("G15").Value = 1200; $ws.Range("H15").Value = $false
$ws.Range("A16").Value = 6; $ws.Range("C16").Value = 50; $ws.Range("G16").Value = 128; $ws.Range("H16").Value = $true
$ws.Range("A17").Value = 10; $ws.Range("C17").Value = 33; $ws.Range("G17").Value = 34; $ws.Range("H17").Value = $false
$ws.Range("A18").Value = 12; $ws.Range("C18").Value = 10; $ws.Range("G18").Value = -6; $ws.Range("H18").Value = $false
$ws.Range("A19").Value = 13; $ws.Range("C19").Value = 4; $ws.Range("G19").Value = 1; $ws.Range("H19").Value = $false
